$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 4 values ---
$ws.Range("A4").Value2 = "Raxmatov islom Xusan o'g'li"
$ws.Range("C4").Value2 = "O'zbek tili"
$ws.Range("E4").Value2 = "AD5538879"

# F4 and I4 are numeric-looking strings; force text storage without leaving
# a lasting number-format style on the cell.
$cellF4 = $ws.Range("F4")
$cellF4.NumberFormat = "@"
$cellF4.Value2 = "31408950420035"
$cellF4.ClearFormats()

$ws.Range("G4").Value2 = "Toshkent viloyati"
$ws.Range("H4").Value2 = "Olmaliq shahri"

$cellI4 = $ws.Range("I4")
$cellI4.NumberFormat = "@"
$cellI4.Value2 = "998931898188"
$cellI4.ClearFormats()

# --- Remove row 5 entirely (shifts dimension from A1:M5 to A1:M4) ---
$ws.Rows.Item(5).Delete()

# --- Update the selection shown in the sheet view ---
$rng = $ws.Range($ws.Rows.Item(4), $ws.Rows.Item(7))
$rng.Select()
